# Iterator, Object, Automation Overview, Framework Overview
#
# The workbook originally has two sheets: "Objects" (sheet1) and "Sheet1"
# (sheet2, currently empty). This script:
#   1. Renames "Sheet1" -> "Maps"
#   2. Populates "Maps" with a remapped copy of the "Objects" data: the
#      same rows/columns, but the "Zip Code" column (C) and "Sub-title2"
#      column (G) are relocated to columns L and N respectively, leaving
#      C, G, K, M empty.
#   3. Re-creates the hyperlinks on column E (About You links) for the
#      new sheet.
#   4. Updates the selections/active sheet so that "Maps" ends up the
#      active (selected) tab with H8-ish selection state matching the
#      target, and "Objects" keeps a non-active selection at H13.

$wb = $excel.ActiveWorkbook
$objects = $wb.Worksheets.Item(1)
$maps = $wb.Worksheets.Item(2)

# --- Rename the second sheet ---
$maps.Name = "Maps"

# --- Copy formatting (borders/fonts/number formats) from Objects!A1:J6 ---
$objects.Range("A1:J6").Copy()
$maps.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

# --- Copy values (and shared-string reuse) from Objects!A1:J6 ---
$objects.Range("A1:J6").Copy()
$maps.Range("A1").PasteSpecial(-4163)   # xlPasteValues

# --- Relocate column C (Zip Code) to column L, and G (Sub-title2) to N ---
$maps.Range("C1:C6").Cut($maps.Range("L1:L6"))
$maps.Range("G1:G6").Cut($maps.Range("N1:N6"))

# The cut leaves the old cells formatted-but-empty; clear them completely
# so they don't appear in the sheet at all (matching the target layout).
$maps.Range("C1:C6").Clear()
$maps.Range("G1:G6").Clear()

# --- Hyperlinks on the "About You" URL column (E2:E6) ---
$maps.Hyperlinks.Add($maps.Range("E2"), "https://sales.geico.com/quote") | Out-Null
$maps.Hyperlinks.Add($maps.Range("E3"), "https://sales.geico.com/quote") | Out-Null
$maps.Hyperlinks.Add($maps.Range("E4"), "https://sales.geico.com/quote") | Out-Null
$maps.Hyperlinks.Add($maps.Range("E5"), "https://sales.geico.com/quote") | Out-Null
$maps.Hyperlinks.Add($maps.Range("E6"), "https://sales.geico.com/quote") | Out-Null

# --- Selections: Objects keeps H13 selected (no longer the active tab);
#     Maps becomes the active tab, selected at K8 ---
$objects.Range("H13").Select()
$maps.Range("K8").Select()
$maps.Activate()
